$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Ipari switch" / wireless / internet-router block (rows 1-12, cols G-M)
$ws.Range("G1").Value = "internet router"
$ws.Range("I1").Value = "wireless kifele"
$ws.Range("J1").Value = "172.0.0.2"
$ws.Range("K1").Value = "255.255.255.252"

$ws.Range("G2").Value = "g0/2"
$ws.Range("I2").Value = "wireless belso"
$ws.Range("J2").Value = "172.16.0.1"
$ws.Range("K2").Value = "255.255.255.192"
$ws.Range("L2").Value = "dhcpvel oszt 10-tol"
$ws.Range("M2").Value = "52 gep tud igy osszesen csatlakozni"

$ws.Range("G3").Value = "172.0.0.5/30"
$ws.Range("I3").Value = "Ipari switch"
$ws.Range("J3").Value = "172.16.0.62"
$ws.Range("K3").Value = "255.255.255.192"
$ws.Range("L3").Value = "admin"

$ws.Range("G4").Value = "s0/3/1"

$ws.Range("G5").Value = "172.0.0.14"
$ws.Range("I5").Value = "wifi"

$ws.Range("G6").Value = "s0/2/0"
$ws.Range("I6").Value = "ssid"
$ws.Range("J6").Value = "ipari_wifi"

$ws.Range("G7").Value = "172.0.0.17"
$ws.Range("I7").Value = "wpa2-psk"
$ws.Range("J7").Value = "cisco123"
$ws.Range("K7").Value = "2.4ghz, 5ghz"

$ws.Range("I10").Value = "tuzfal"

$ws.Range("I11").Value = "kifele g1/1"
$ws.Range("J11").Value = "172.0.0.6"
$ws.Range("K11").Value = "255.255.255.252"

$ws.Range("I12").Value = "befele g1/2"
$ws.Range("J12").Value = "172.0.0.1"
$ws.Range("K12").Value = "255.255.255.252"

# New "Iroda1 router" / "Master router" / "Backup router" block (rows 15-26, cols G-L)
$ws.Range("G15").Value = "Iroda1 router"

$ws.Range("G16").Value = "g0/0"
$ws.Range("H16").Value = "ipv6"
$ws.Range("I16").Value = "2001:db8:1::1/64"
$ws.Range("J16").Value = "s0/3/0"
$ws.Range("K16").Value = "ip"
$ws.Range("L16").Value = "172.0.0.13"

$ws.Range("G17").Value = "Master router"

$ws.Range("G18").Value = "virtual cím"
$ws.Range("I18").Value = "FE80::10"

$ws.Range("G19").Value = "priority"
$ws.Range("I19").Value = 120

$ws.Range("G20").Value = "preempt"

$ws.Range("G21").Value = "iroda2 router"
$ws.Range("J21").Value = "s0/3/0"
$ws.Range("K21").Value = "ip"
$ws.Range("L21").Value = "172.0.0.18"

$ws.Range("G22").Value = "g0/0"
$ws.Range("H22").Value = "ipv6"
$ws.Range("I22").Value = "2001:db8:1::2/64"

$ws.Range("G23").Value = "Backup router"

$ws.Range("G24").Value = "virtual cím"
$ws.Range("I24").Value = "FE80::10"

$ws.Range("G25").Value = "priority"
$ws.Range("I25").Value = 100

$ws.Range("G26").Value = "preempt"

# Match the author's final selection state
$ws.Range("I17").Select()
